$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value  = -7.904500000000008
$ws.Range("A9").Value  = -22.0605
$ws.Range("D9").Value  = -8.3499
$ws.Range("A18").Value = -22.19380000000001
$ws.Range("A20").Value = -21.58769999999998
$ws.Range("D23").Value = -8.0336
$ws.Range("D24").Value = -7.202299999999997
$ws.Range("D26").Value = -7.866099999999999
$ws.Range("A27").Value = -22.02900000000001
$ws.Range("D34").Value = -7.886300000000003
$ws.Range("D35").Value = -7.737600000000001
$ws.Range("D48").Value = -7.542300000000001
$ws.Range("D52").Value = -7.881400000000006
$ws.Range("D66").Value = -7.196800000000001
$ws.Range("D67").Value = -7.591599999999999
$ws.Range("A69").Value = -21.58359999999998
$ws.Range("A76").Value = -19.61739999999998
$ws.Range("D80").Value = -7.921300000000003
$ws.Range("A82").Value = -21.83349999999999
$ws.Range("D99").Value = -8.0844
